$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format before writing, so numeric-looking
# strings (e.g. "0.9996", "238.80") are stored as text, matching the
# original inlineStr cell type instead of being auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '30.909.98'
$ws.Range("E2").Value = '  +0.60%  '

$ws.Range("D3").Value = '1.907.79'
$ws.Range("E3").Value = '  +0.71%  '

$ws.Range("D4").Value = '0.9996'
$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").Value = '238.80'
$ws.Range("E5").Value = '  -3.26%  '

$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  +0.02%  '

$ws.Range("E7").Value = '  -0.56%  '

$ws.Range("D8").Value = '0.2968'
$ws.Range("E8").Value = '  +0.47%  '

$ws.Range("D9").Value = '0.06760'
$ws.Range("E9").Value = '  -0.73%  '

$ws.Range("D10").Value = '1.920.82'
$ws.Range("E10").Value = '  +1.40%  '

$ws.Range("D11").Value = '17.07'
$ws.Range("E11").Value = '  -1.25%  '

$ws.Range("D12").Value = '0.07296'
$ws.Range("E12").Value = '  +0.55%  '

$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '5.135'
$ws.Range("E13").Value = '  +0.62%  '

$ws.Range("B14").Value = 'Litecoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D14").Value = '89.79'
$ws.Range("E14").Value = '  -2.56%  '

$ws.Range("D15").Value = '0.6702'
$ws.Range("E15").Value = '  -1.79%  '

$ws.Range("D16").Value = '30.876.09'
$ws.Range("E16").Value = '  +0.57%  '

$ws.Range("D17").Value = '0.000007949'
$ws.Range("E17").Value = '  -0.58%  '

$ws.Range("D18").Value = '13.44'
$ws.Range("E18").Value = '  +0.89%  '

$ws.Range("D19").Value = '0.9996'
$ws.Range("E19").Value = '  -0.09%  '

$ws.Range("D20").Value = '2.193.68'
$ws.Range("E20").Value = '  +2.54%  '

$ws.Range("D21").Value = '0.9989'
$ws.Range("E21").Value = '  -0.25%  '

$ws.Range("D22").Value = '5.105'
$ws.Range("E22").Value = '  +5.21%  '

$ws.Range("D23").Value = '206.81'
$ws.Range("E23").Value = '  +7.30%  '

$ws.Range("E24").Value = '  +2.26%  '

$ws.Range("D25").Value = '9.642'
$ws.Range("E25").Value = '  +2.43%  '

$ws.Range("D26").Value = '157.54'

$ws.Range("D27").Value = '18.85'
$ws.Range("E27").Value = '  -1.73%  '

$ws.Range("D28").Value = '1.968'
$ws.Range("E28").Value = '  +2.10%  '

$ws.Range("D29").Value = '1.424'
$ws.Range("E29").Value = '  +1.71%  '

$ws.Range("D30").Value = '4.337'
$ws.Range("E30").Value = '  -0.79%  '

$ws.Range("D31").Value = '0.09162'
$ws.Range("E31").Value = '  +1.61%  '

$ws.Range("D32").Value = '4.054'
$ws.Range("E32").Value = '  +0.23%  '

$ws.Range("D33").Value = '0.05178'
$ws.Range("E33").Value = '  -0.55%  '

$ws.Range("D34").Value = '0.7531'
$ws.Range("E34").Value = '  +0.97%  '

$ws.Range("E35").Value = '  -1.05%  '

$ws.Range("D36").Value = '2.708'
$ws.Range("E36").Value = '  -0.50%  '

$ws.Range("D37").Value = '0.01842'
$ws.Range("E37").Value = '  -0.97%  '

$ws.Range("D38").Value = '2.723'
$ws.Range("E38").Value = '  +1.74%  '

$ws.Range("D39").Value = '0.9258'
$ws.Range("E39").Value = '  -1.27%  '

$ws.Range("D40").Value = '2.106'
$ws.Range("E40").Value = '  -2.57%  '

$ws.Range("D41").Value = '0.4493'
$ws.Range("E41").Value = '  +1.13%  '

$ws.Range("D42").Value = '107.07'
$ws.Range("E42").Value = '  +0.49%  '

$ws.Range("D43").Value = '5.869'
$ws.Range("E43").Value = '  +1.35%  '

$ws.Range("E44").Value = '  +0.57%  '

$ws.Range("D45").Value = '7.737'
$ws.Range("E45").Value = '  +0.55%  '

$ws.Range("D46").Value = '0.1377'
$ws.Range("E46").Value = '  +2.63%  '

$ws.Range("D47").Value = '67.65'
$ws.Range("E47").Value = '  +16.53%  '

$ws.Range("B48").Value = 'Elrond'
$ws.Range("C48").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D48").Value = '35.26'
$ws.Range("E48").Value = '  +4.80%  '

$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '8.989'
$ws.Range("E49").Value = '  +2.41%  '

$ws.Range("B50").Value = 'Decentraland'
$ws.Range("C50").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D50").Value = '0.4099'
$ws.Range("E50").Value = '  +3.67%  '

$ws.Range("E51").Value = '  +0.75%  '

# Restore column D cells to the default (unstyled) look so no spurious
# style index is left behind on the cells after the Text-format trick.
$ws.Range("D2:D51").Style = "Normal"
